$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.740.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.98%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.302.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.08%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.92"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.20"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.29%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.300.65"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.73%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.71%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.56%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.68"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.846.66"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.14%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.61%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.296.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.827.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.92%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.54%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.32%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.75%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.18%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.68"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.42"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.50%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.13%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.12%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.81%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.12"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.91%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.69"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.97%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.55"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.97%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.15%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.37"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.12%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0742"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.50%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0401"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "432.49"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.074.83"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.56%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.35"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.84%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.29%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.114"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.265"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.38%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.32"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.59%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.30"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.62%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.02"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +12.03%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.115"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.60%  "
